# Implement a new variant of relative formula type ("RelativeF Samples")
# that only accepts formula cells (not hardcoded constants).
#
# This inserts a new worksheet, "RelativeF Samples", right before the
# existing "Relative Samples" sheet. The new sheet reuses the same
# data shape as "Relative Samples" but column A now holds formulas
# (referencing new columns E/F) instead of hardcoded constants.

$wb = $excel.ActiveWorkbook

$relative = $wb.Worksheets.Item("Relative Samples")

# Insert the new worksheet immediately before "Relative Samples".
$relativeF = $wb.Worksheets.Add($relative)
$relativeF.Name = "RelativeF Samples"

# Row 2
$relativeF.Range("A2").Formula = "=E2-F2"
$relativeF.Range("B2").Value = 1001
$relativeF.Range("C2").Value = 2002
$relativeF.Range("E2").Value = 4004
$relativeF.Range("F2").Value = 1001

# Row 3
$relativeF.Range("A3").Formula = "=E3+F3"
$relativeF.Range("B3").Value = 2002
$relativeF.Range("C3").Value = 2002
$relativeF.Range("E3").Value = 1000
$relativeF.Range("F3").Value = 3005

# Row 4
$relativeF.Range("A4").Formula = "=E4+F4"
$relativeF.Range("B4").Value = 1000
$relativeF.Range("C4").Value = 2000
$relativeF.Range("E4").Value = 1000
$relativeF.Range("F4").Value = 2005

# Row 5
$relativeF.Range("A5").Formula = "=E5+F5"
$relativeF.Range("B5").Value = 1000
$relativeF.Range("C5").Value = 2000
$relativeF.Range("E5").Value = 4444
$relativeF.Range("F5").Value = 4444

# Row 6 (hardcoded "not_ok" marker row, no formula/E/F values)
$relativeF.Range("A6").Value = 2000
$relativeF.Range("B6").Value = "not_ok"
$relativeF.Range("C6").Value = 1000
$relativeF.Range("D6").Value = 2000

# Re-fetch "Relative Samples" by name: after Worksheets.Add() shifted its
# index, the original $relative handle no longer tracks the right sheet.
$relative = $wb.Worksheets.Item("Relative Samples")

# Update the selection/tab state on the (now relocated) "Relative Samples"
# sheet first so it is no longer the active/selected tab.
$relative.Activate()
$relative.Range("G8").Select()

# Make the new sheet the active tab with A2 selected.
$relativeF.Activate()
$relativeF.Range("A2").Select()
